# Add a "Difficulty" column (C) next to the existing Letters/Done columns,
# rating the transcription difficulty of each letter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled like the other header cells (light blue fill).
$ws.Range("C1").Value = "Difficulty"
$ws.Range("C1").Interior.Color = 15652797

# Per-row difficulty ratings (rows already marked "Done" in column B are left blank).
$ws.Range("C2").Value = "high"       # liv_000493
$ws.Range("C3").Value = "high"       # liv_000586
$ws.Range("C6").Value = "low-med"    # liv_001342
$ws.Range("C8").Value = "low-med"    # liv_001509
$ws.Range("C9").Value = "high"       # liv_001572
$ws.Range("C13").Value = "low"       # liv_002698
$ws.Range("C10").Value = "med-high"  # liv_001965
$ws.Range("C11").Value = "low-med"   # liv_002208
$ws.Range("C12").Value = "high"      # liv_002697
$ws.Range("C14").Value = "low"       # liv_002699
$ws.Range("C16").Value = "low-med"   # liv_002704
$ws.Range("C17").Value = "low-med"   # liv_002705
$ws.Range("C18").Value = "medium"    # liv_002706
$ws.Range("C19").Value = "low-med"   # liv_002707

$ws.Range("C4").Select()
